$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.570.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5067"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3952"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09806"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +5.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.547"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.922.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.583"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06665"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.325"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.634.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.741"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.133.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.56%  "
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +7.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.761"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.892"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06807"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02446"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6444"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.191"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.32%  "
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.820"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.51%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  +3.36%  "
